$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column O (15) to make room for the new
#    "Age Match?" column. This shifts the existing "Remarks" (O) and
#    "References" (P) headers/columns one place to the right (-> P, Q)
#    and widens the M1:P1 "Author Review" merge to M1:Q1 automatically.
$ws.Columns.Item(15).Insert()

# 2. Give the new column its header text.
$ws.Cells.Item(2, 15).Value2 = "Age Match?"

# 3. Update the "Extent" column (F) text for every row that currently
#    reads "oh, " so that it reads the full state list.
$ohRows = @(3,4,5,6,7,8,9,11,20,21,22,27,32,33,34,43,44,46,48,49,50,52,53,54,57,58,62,63,64,68,69,70,71,72,74,76)
foreach ($r in $ohRows) {
    $ws.Cells.Item($r, 6).Value2 = "OH, PA, WV, VA"
}

# 4. Flip the "Extent Match?" column (M) answer from "no" to "yes" for
#    the rows that were re-reviewed.
$yesRows = @(9,10,11,12,13,17,18,35,36)
foreach ($r in $yesRows) {
    $ws.Cells.Item($r, 13).Value2 = "yes"
}

# 5. Row 54's Age cell (E54) switches from a numeric 15 to the text "15".
$ageCell = $ws.Cells.Item(54, 5)
$ageCell.NumberFormat = "@"
$ageCell.Value2 = "15"
